$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14197
$ws.Range("F4").Value = 14197
$ws.Range("F5").Value = 14303
$ws.Range("F10").Value = 573
$ws.Range("F14").Value = 1542
$ws.Range("F18").Value = 1831
$ws.Range("F22").Value = 564
$ws.Range("F24").Value = 3320
$ws.Range("F33").Value = 1391
$ws.Range("F34").Value = 101
$ws.Range("F35").Value = 148
$ws.Range("F36").Value = 4812
$ws.Range("F37").Value = 4849
$ws.Range("F42").Value = 3289
$ws.Range("F45").Value = 339
$ws.Range("F48").Value = 4422

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 119

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7560
$ws.Range("F4").Value = 769

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7560
$ws.Range("F5").Value = 769
$ws.Range("F6").Value = 14197
$ws.Range("F7").Value = 14303
$ws.Range("F12").Value = 119
$ws.Range("F15").Value = 1542
$ws.Range("F18").Value = 1831
$ws.Range("F21").Value = 3320
$ws.Range("F32").Value = 1391
$ws.Range("F33").Value = 101
$ws.Range("F34").Value = 4812
$ws.Range("F35").Value = 4849
$ws.Range("F40").Value = 3289
$ws.Range("F42").Value = 339
$ws.Range("F46").Value = 4422
